$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Aksobhya: update image filename
$ws.Range("F2").Value = "BAkso008.JPG"

# Row 5 - was "Thangkas  - groß" / "001.JPG" -> "Thangkas - Big Tsagil" / multiple images
$ws.Range("B5").Value = "Thangkas - Big Tsagil"

# Row 4 - was "Malas  - klein" / "001.JPG" -> "Malas - Arm Mala" / "MAM020.JPG"
$ws.Range("B4").Value = "Malas - Arm Mala"
$ws.Range("F4").Value = "MAM020.JPG"

# Row 3 - was "Buddhas - Manjusri" / "002.jpg" -> "Buddhas - Shakyamuni" / multiple images
$ws.Range("B3").Value = "Buddhas - Shakyamuni"
$ws.Range("F3").Value = "BSha001.JPG, BSha002.JPG, BSha003.JPG"

$ws.Range("F5").Value = "TBT020.JPG, TBT019.JPG, TBT018.JPG"

# F5 loses its explicit "Text" number-format style in the target workbook
$ws.Range("A5").Copy()
$ws.Range("F5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the active selection to match the edited workbook
$ws.Range("G14").Select()
